# edit.ps1 - reproduces:
#   1) table style swap on the three data tables (slides 14, 15, 16):
#        {E8BCD833-53B3-4809-AA01-6D16EE31D932} -> {1A280DE8-6E74-4BA9-8F39-36E026ECAECF}
#   2) the deck's theme colour scheme reverting from the "Integral" / "Red Violet"
#      palette to the default "Office" palette.

$p = $ppt.ActivePresentation

# --- 1. Re-style the three tables -----------------------------------------
$newStyleId = "{1A280DE8-6E74-4BA9-8F39-36E026ECAECF}"

foreach ($slideIdx in 14, 15, 16) {
    $slide = $p.Slides.Item($slideIdx)
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shape = $slide.Shapes.Item($i)
        if ($shape.HasTable) {
            $shape.Table.ApplyStyle($newStyleId)
        }
    }
}

# --- 2. Restore the default Office colour scheme on the deck theme --------
# (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink) as packed BGR integers
# suitable for the VBA-style ColorScheme.Colors(n).RGB setter.
$officeColorsHex = @(
    "000000",  # dk1
    "FFFFFF",  # lt1
    "44546A",  # dk2
    "E7E6E6",  # lt2
    "5B9BD5",  # accent1
    "ED7D31",  # accent2
    "A5A5A5",  # accent3
    "FFC000",  # accent4
    "4472C4",  # accent5
    "70AD47",  # accent6
    "0563C1",  # hlink
    "954F72"   # folHlink
)

$cs = $p.SlideMaster.ColorScheme
for ($i = 1; $i -le $officeColorsHex.Count; $i++) {
    $hex = $officeColorsHex[$i - 1]
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    $bgr = $r + ($g * 256) + ($b * 65536)
    $cs.Colors($i).RGB = $bgr
}
